# The edit cyclically rotates the text of several paragraphs:
#   Objetivos             <- old Programa resumido
#   Docente(s)             <- old Objetivos
#   Programa resumido      <- old Programa
#   Programa                <- old Metodo
#   Metodo                   <- old Criterio
#   Criterio                  <- old Norma de recuperacao
#   Norma de recuperacao    <- old Bibliografia
#   Bibliografia              <- old Docente(s)
#
# Because several of these texts would otherwise overwrite one another
# (the source of one replacement is the target of another), first swap
# every original text out for a unique placeholder token, and only then
# replace each placeholder with its real destination text.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Step 1: move every original value into a unique placeholder so that
# the subsequent replacements do not collide with each other.
Replace-Text "Avaliar casos de impacto ambiental que proporcionaram singularidades em algumas de suas etapas." "__PLACEHOLDER_0__"
Replace-Text "5840938 - Marcelo Rodrigues de Holanda" "__PLACEHOLDER_1__"
Replace-Text "Estudos de caso: avaliar a singularidade do caso apresentado, com a necessária identificação das etapas e das peculiaridades que os fizeram próprios a serem aplicados em estudo de caso." "__PLACEHOLDER_2__"
Replace-Text "Estudos de casos específicos e as etapas necessárias na avaliação de um impacto ambiental." "__PLACEHOLDER_3__"
Replace-Text "Aula expositiva e exercícios dirigidos." "__PLACEHOLDER_4__"
Replace-Text "Média ponderada de exercícios e provas." "__PLACEHOLDER_5__"
Replace-Text "Prova única com nota igual ou superior a 5,0." "__PLACEHOLDER_6__"
Replace-Text "Estudos de caso: EPIA de origem." "__PLACEHOLDER_7__"

# Step 2: place each paragraph's final text value.
Replace-Text "__PLACEHOLDER_0__" "Estudos de caso: avaliar a singularidade do caso apresentado, com a necessária identificação das etapas e das peculiaridades que os fizeram próprios a serem aplicados em estudo de caso."
Replace-Text "__PLACEHOLDER_1__" "Avaliar casos de impacto ambiental que proporcionaram singularidades em algumas de suas etapas."
Replace-Text "__PLACEHOLDER_2__" "Estudos de casos específicos e as etapas necessárias na avaliação de um impacto ambiental."
Replace-Text "__PLACEHOLDER_3__" "Aula expositiva e exercícios dirigidos."
Replace-Text "__PLACEHOLDER_4__" "Média ponderada de exercícios e provas."
Replace-Text "__PLACEHOLDER_5__" "Prova única com nota igual ou superior a 5,0."
Replace-Text "__PLACEHOLDER_6__" "Estudos de caso: EPIA de origem."
Replace-Text "__PLACEHOLDER_7__" "5840938 - Marcelo Rodrigues de Holanda"

Write-Output "done"
